$d = $word.ActiveDocument

# Locate the "Navid Shekoufa" contribution line: "...201577806,<tab>12/40"
# and fix the wrongly entered contribution score "12/40" -> "13/40" by
# changing just the "2" digit run that sits right before the _GoBack
# bookmark (so the bookmark and surrounding runs stay untouched).
$rng = $d.Content
$found = $rng.Find.Execute("201577806,`t1")
if ($found) {
    $target = $d.Range($rng.End, $rng.End + 1)
    if ($target.Text -eq "2") {
        $target.Text = "3"
    }
}
